# Refresh the "cryptos" price/volume snapshot (GitHub Actions scheduled update).
# Price cells in column D are stored as plain text (e.g. "64.079.27",
# "14.30", "0.110") so that thousand-separator dots and trailing zeros
# survive round-tripping. Writing a numeric-looking string straight into
# `.Value` would get auto-coerced to a real number (losing the trailing
# zero / European grouping dots), so for any D-column value that parses as
# a number we first force the cell to Text format, assign the literal
# string, then snap the cell style back to "Normal" so no stray formatting
# is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.079.27"
$ws.Range("E2").Value = "  +1.61%  "
$ws.Range("D3").Value = "3.309.05"
$ws.Range("E3").Value = "  +5.86%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.91%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "3.304.01"
$ws.Range("E8").Value = "  +5.80%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.522"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.45%  "
$ws.Range("E10").Value = "  +2.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.48"
$ws.Range("D11").Style = "Normal"
$ws.Range("E12").Value = "  +2.62%  "
$ws.Range("E13").Value = "  -0.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.96"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.82%  "
$ws.Range("D15").Value = "3.850.48"
$ws.Range("E16").Value = "  +1.18%  "
$ws.Range("D17").Value = "3.306.30"
$ws.Range("E17").Value = "  +5.92%  "
$ws.Range("D18").Value = "64.147.96"
$ws.Range("E18").Value = "  +1.68%  "
$ws.Range("E19").Value = "  +2.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "483.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.96%  "
$ws.Range("E22").Value = "  +6.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.03"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.54"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.49"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.36%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("E27").Value = "  +2.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.36"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.30"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.55%  "
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.15"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.75"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.88%  "
$ws.Range("E33").Value = "  -1.70%  "
$ws.Range("E34").Value = "  +1.09%  "
$ws.Range("E35").Value = "  +1.96%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "53.40"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.72%  "
$ws.Range("E38").Value = "  +3.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0401"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "431.57"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.64%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "3.013.71"
$ws.Range("E41").Value = "  +4.55%  "
$ws.Range("B42").Value = "Cosmos"
$ws.Range("C42").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.46"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.10%  "
$ws.Range("E43").Value = "  +2.74%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.110"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.270"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.38%  "
$ws.Range("E46").Value = "  +4.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "26.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.36%  "
$ws.Range("E48").Value = "  +0.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.69"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +14.67%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.115"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.52%  "
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.33"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.12%  "
